# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "57.764.39"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D2:E2").Style = "Normal"

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "2.348.94"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D3:E3").Style = "Normal"

# Row 4
$ws.Range("E4:E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E4:E4").Style = "Normal"

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "547.65"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D5:E5").Style = "Normal"

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "132.02"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("D6:E6").Style = "Normal"

# Row 7
$ws.Range("E7:E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E7:E7").Style = "Normal"

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D8:E8").Style = "Normal"

# Row 9
$ws.Range("E9:E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("E9:E9").Style = "Normal"

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "5.61"
$ws.Range("E10").Value = "  +3.98%  "
$ws.Range("D10:E10").Style = "Normal"

# Row 11
$ws.Range("E11:E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E11:E11").Style = "Normal"

# Row 12
$ws.Range("E12:E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("E12:E12").Style = "Normal"

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "23.80"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D13:E13").Style = "Normal"

# Row 14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "2.766.16"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D14:E14").Style = "Normal"

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "57.691.91"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D15:E15").Style = "Normal"

# Row 16
$ws.Range("E16:E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E16:E16").Style = "Normal"

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "2.387.75"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D17:E17").Style = "Normal"

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "10.96"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D18:E18").Style = "Normal"

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "4.28"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D19:E19").Style = "Normal"

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "329.31"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D20:E20").Style = "Normal"

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D21:E21").Style = "Normal"

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D22:E22").Style = "Normal"

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "63.71"
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("D23:E23").Style = "Normal"

# Row 24
$ws.Range("E24:E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E24:E24").Style = "Normal"

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D25:E25").Style = "Normal"

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "8.22"
$ws.Range("E26").Value = "  -3.41%  "
$ws.Range("D26:E26").Style = "Normal"

# Row 27
$ws.Range("E27:E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.43%  "
$ws.Range("E27:E27").Style = "Normal"

# Row 28
$ws.Range("E28:E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("E28:E28").Style = "Normal"

# Row 29
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "171.17"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D29:E29").Style = "Normal"

# Row 30
$ws.Range("E30:E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E30:E30").Style = "Normal"

# Row 31
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D31:E31").Style = "Normal"

# Row 32
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "18.36"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D32:E32").Style = "Normal"

# Row 33
$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D33:E33").Style = "Normal"

# Row 34
$ws.Range("E34:E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E34:E34").Style = "Normal"

# Row 35
$ws.Range("E35:E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("E35:E35").Style = "Normal"

# Row 36
$ws.Range("E36:E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("E36:E36").Style = "Normal"

# Row 37
$ws.Range("E37:E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +16.94%  "
$ws.Range("E37:E37").Style = "Normal"

# Row 38
$ws.Range("E38:E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E38:E38").Style = "Normal"

# Row 39
$ws.Range("B39:E39").NumberFormat = "@"
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "40.35"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("B39:E39").Style = "Normal"

# Row 40
$ws.Range("B40:E40").NumberFormat = "@"
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.59"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("B40:E40").Style = "Normal"

# Row 41
$ws.Range("B41:E41").NumberFormat = "@"
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "141.25"
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("B41:E41").Style = "Normal"

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "3.63"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D42:E42").Style = "Normal"

# Row 43
$ws.Range("B43:E43").NumberFormat = "@"
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "287.33"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("B43:E43").Style = "Normal"

# Row 44
$ws.Range("B44:E44").NumberFormat = "@"
$ws.Range("B44").Value = "Polygon"
$ws.Range("C44").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D44").Value = "0.426"
$ws.Range("E44").Value = "  +11.74%  "
$ws.Range("B44:E44").Style = "Normal"

# Row 45
$ws.Range("B45:E45").NumberFormat = "@"
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0951"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("B45:E45").Style = "Normal"

# Row 46
$ws.Range("B46:E46").NumberFormat = "@"
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0511"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("B46:E46").Style = "Normal"

# Row 47
$ws.Range("B47:E47").NumberFormat = "@"
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.566"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("B47:E47").Style = "Normal"

# Row 48
$ws.Range("B48:E48").NumberFormat = "@"
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "18.65"
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("B48:E48").Style = "Normal"

# Row 49
$ws.Range("B49:E49").NumberFormat = "@"
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0220"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("B49:E49").Style = "Normal"

# Row 50
$ws.Range("B50:E50").NumberFormat = "@"
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "11.09"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B50:E50").Style = "Normal"

# Row 51
$ws.Range("B51:E51").NumberFormat = "@"
$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").Value = "4.70"
$ws.Range("E51").Value = "  -0.08%  "
$ws.Range("B51:E51").Style = "Normal"
